$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand-expressing cell counts / rates / expression values for the
# "ECs" sending cluster (rows 2-4), reflecting the new TPM input.
$ws.Range("E2:E4").Value = 1
$ws.Range("F2:F4").Value = 0.3333333333333333
$ws.Range("G2:G4").Value = 0.004531666666666666
$ws.Range("H2:H4").Value = 0.013595

# Updated receptor average / total expression values for the "ECs" target
# cluster (rows 2, 5, 8).
$ws.Range("M2").Value = 0.366183
$ws.Range("N2").Value = 1.098549
$ws.Range("M5").Value = 0.366183
$ws.Range("N5").Value = 1.098549
$ws.Range("M8").Value = 0.366183
$ws.Range("N8").Value = 1.098549

# Recomputed ligand-derived specificity (I, J) for every row, since the
# total ligand expression across sending clusters changed.
$ws.Range("I2").Value = 0.00937711019466729
$ws.Range("J2").Value = 0.00937711019466729
$ws.Range("I3").Value = 0.00937711019466729
$ws.Range("J3").Value = 0.00937711019466729
$ws.Range("I4").Value = 0.00937711019466729
$ws.Range("J4").Value = 0.00937711019466729
$ws.Range("I5").Value = 0.2767209704464111
$ws.Range("J5").Value = 0.2767209704464111
$ws.Range("I6").Value = 0.2767209704464111
$ws.Range("J6").Value = 0.2767209704464111
$ws.Range("I7").Value = 0.2767209704464111
$ws.Range("J7").Value = 0.2767209704464111
$ws.Range("I8").Value = 0.7139019193589216
$ws.Range("J8").Value = 0.7139019193589216
$ws.Range("I9").Value = 0.7139019193589216
$ws.Range("J9").Value = 0.7139019193589216
$ws.Range("I10").Value = 0.7139019193589216
$ws.Range("J10").Value = 0.7139019193589216

# Recomputed receptor-derived specificity (O, P) for every row, since the
# total receptor expression across target clusters changed.
$ws.Range("O2").Value = 0.0639836884691917
$ws.Range("P2").Value = 0.0639836884691917
$ws.Range("O3").Value = 0.2777364052521014
$ws.Range("P3").Value = 0.2777364052521014
$ws.Range("O4").Value = 0.6582799062787069
$ws.Range("P4").Value = 0.6582799062787069
$ws.Range("O5").Value = 0.0639836884691917
$ws.Range("P5").Value = 0.0639836884691917
$ws.Range("O6").Value = 0.2777364052521014
$ws.Range("P6").Value = 0.2777364052521014
$ws.Range("O7").Value = 0.6582799062787069
$ws.Range("P7").Value = 0.6582799062787069
$ws.Range("O8").Value = 0.0639836884691917
$ws.Range("P8").Value = 0.0639836884691917
$ws.Range("O9").Value = 0.2777364052521014
$ws.Range("P9").Value = 0.2777364052521014
$ws.Range("O10").Value = 0.6582799062787069
$ws.Range("P10").Value = 0.6582799062787069

# Recomputed edge expression weights (Q, R = ligand * receptor expression)
$ws.Range("Q2").Value = 0.001659419295
$ws.Range("R2").Value = 0.014934773655
$ws.Range("Q3").Value = 0.007203103803888889
$ws.Range("R3").Value = 0.06482793423500001
$ws.Range("Q4").Value = 0.01707251338777778
$ws.Range("R4").Value = 0.15365262049
$ws.Range("Q5").Value = 0.048969896712
$ws.Range("R5").Value = 0.440729070408
$ws.Range("Q8").Value = 0.12633557622
$ws.Range("R8").Value = 1.13702018598

# Recomputed edge derived specificity (S, T = ligand specificity * receptor
# specificity) for every row.
$ws.Range("S2").Value = 0.0005999820974368734
$ws.Range("T2").Value = 0.0005999820974368734
$ws.Range("S3").Value = 0.002604364877119726
$ws.Range("T3").Value = 0.002604364877119726
$ws.Range("S4").Value = 0.006172763220110691
$ws.Range("T4").Value = 0.006172763220110691
$ws.Range("S5").Value = 0.01770562836593557
$ws.Range("T5").Value = 0.01770562836593557
$ws.Range("S6").Value = 0.07685548758965922
$ws.Range("T6").Value = 0.07685548758965922
$ws.Range("S7").Value = 0.1821598544908163
$ws.Range("T7").Value = 0.1821598544908163
$ws.Range("S8").Value = 0.04567807800581925
$ws.Range("T8").Value = 0.04567807800581925
$ws.Range("S9").Value = 0.1982765527853225
$ws.Range("T9").Value = 0.1982765527853225
$ws.Range("S10").Value = 0.4699472885677798
$ws.Range("T10").Value = 0.4699472885677798
